$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 24, pushing the existing rows 24..80 down to 25..81.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44557
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 375
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 10667
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región de La Araucanía"
$ws.Range("P24").Value = 427
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
